$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 label
$ws.Range("A2").Value = "id_DK_Central_EP"

# Update row 3 label and value (aggregated)
$ws.Range("A3").Value = "id_DK_Decentral_EP"
$ws.Range("B3").Value = -2.592592592592593

# Remove row 4 entirely (was id_DK1_SmallDecentral_EP), data now aggregated into row 3
$ws.Rows.Item(4).Delete()
